$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D (Price) and E (Volume/1h) hold numeric-looking text
# (e.g. "306.91", "-0.47%"). The source workbook stores these as literal
# text (inline strings), so force the Text number format before writing
# the value to avoid Excel auto-converting them into real numbers/percentages.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "306.91"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.47%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "38.87"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.81%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.097"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.69%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.50%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.942"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-4.55%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.183"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.61%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.968"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.20%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9312"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.54%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1475"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.27%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1937"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.33%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09181"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.43%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03504"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.39%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09780"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.32%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001395"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.00%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005877"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.79%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.788"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.25%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.458"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.83%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.04%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1303"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.78%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.562"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-5.13%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04362"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.34%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.26%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004282"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-12.94%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.13%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02042"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.36%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05082"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.45%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007526"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.71%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1349"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.77%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002121"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.34%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009082"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-6.69%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006184"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.03%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003103"
